$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($ws1)
$newSheet.Name = "UnitTest"
$newSheet.Range("F4").Value = 43115
$newSheet.Range("F4").Interior.Color = 65535
$newSheet.Range("F4").NumberFormat = "m/d/yy"
